# Import fund and folio amounts
# Inserts a new "Amount (Fund Currency)*" column before the old "Amount *"
# column (which is renamed "Amount (Folio Currency)*"), shifting every
# column from F..L one place to the right (G..M), and adds a new comment
# explaining the new column.

$excel.UserName = "thimm"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Relocate the legacy cell comments that live on row 1 of the
#        columns that are about to shift right, starting from the
#        right-most column so we never clobber a comment we still need
#        to read.
$moves = @(
    @{ From = "L1"; To = "M1" },
    @{ From = "K1"; To = "L1" },
    @{ From = "J1"; To = "K1" },
    @{ From = "I1"; To = "J1" },
    @{ From = "H1"; To = "I1" },
    @{ From = "G1"; To = "H1" },
    @{ From = "F1"; To = "G1" }
)

foreach ($m in $moves) {
    $fromRange = $ws.Range($m.From)
    $text = $fromRange.Comment.Text()
    $fromRange.Comment.Delete()
    $toRange = $ws.Range($m.To)
    $toRange.ClearComments()
    $toRange.AddComment($text) | Out-Null
}

# --- 2. Insert the new column in front of the (old) column F. Excel
#        shifts the column widths, data validation and cell data of
#        F:L to G:M automatically.
$ws.Columns("F:F").Insert()

# --- 3. Rename the old "Amount *" header and add the new header.
$ws.Range("D1").Value = "Amount (Folio Currency)*"
$ws.Range("F1").Value = "Amount (Fund Currency)*"

# --- 4. Match the manually-resized column width that the author gave
#        the new column.
$ws.Columns("F:F").ColumnWidth = 8.4375

# --- 5. Add the new comment that documents the new column.
$ws.Range("F1").ClearComments()
$ws.Range("F1").AddComment("thimm:" + [char]10 + "Amount in fund currency, leave blank if you want the platform to convert the amount based on exchange rates setup on the platform") | Out-Null

# --- 6. Restore the author's final selection.
$ws.Range("F2").Select()
